# Update the YTD CVD value from 0.0639 to 0.0776 across all facility
# sheets that carry it (rows 4/5 for most sheets, rows 2/3 for the
# sheets whose CVD table starts two rows earlier).
$wb = $excel.ActiveWorkbook

$sheetsRows45 = @(
    "Cassville Missouri",
    "Ciserano Italy",
    "Faridabad India",
    "Jiaxing China",
    "Piedras Negras Fasco Mexico",
    "Reynosa Mexico",
    "Sao Paulo Brazil"
)

$sheetsRows23 = @(
    "Mcallen Texas",
    "Piedras Negras Jakel Mexico",
    "Reynosa II"
)

foreach ($name in $sheetsRows45) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E4").Value = 0.0776
    $ws.Range("E5").Value = 0.0776
}

foreach ($name in $sheetsRows23) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E2").Value = 0.0776
    $ws.Range("E3").Value = 0.0776
}
